# Populate the "Eingang" (incoming goods) worksheet with the article table:
# Artikelnummer / Menge / Status header row plus six article rows (AS2005..AS2010).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$data = @(
    @("Artikelnummer", "Menge", "Status"),
    @("AS2005", 7, 0),
    @("AS2006", 8, 0),
    @("AS2007", 9, 0),
    @("AS2008", 10, 1),
    @("AS2009", 11, 0),
    @("AS2010", 12, 0)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

# Matches the saved selection in the edited file (active cell C5).
$ws.Range("C5").Select()
